# New API Query - 2023 Included
# API query to UN performed 11/26/2023. Query modified to include 2023 data.
#
# The underlying "short-url" and "oip" (null -> "-") values changed for
# every data row, and the "hst" column (previously a hard-coded 0) now
# mirrors the "oip" column ("-").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# short-url column (B) changed for all data rows (2-7)
$ws.Range("B2:B7").Value = "vtz2SQ"

# oip column (U): "null" -> "-" for all data rows (2-7)
$ws.Range("U2:U7").Value = "-"

# hst column (V): now also "-", left-aligned like the other text columns
$ws.Range("V2:V7").Value = "-"
$ws.Range("V2:V7").HorizontalAlignment = -4131
